$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Similarity Score"
$ws.Range("B1").Value = "Overhead Reduction"
$ws.Range("C1").Value = "Budget"
